$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Boolean" sheet: split the single trans/BVTQaZ/BVTQaZ.csv and
# trans/VTQaZ/VTQaZ.csv entries into six per-vehicle-type CSV files each.
# ---------------------------------------------------------------------
$wsBool = $wb.Worksheets.Item("Boolean")

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv" -- insert 5 more rows
# below it so there's room for the six BVTQaZ-* entries.
$wsBool.Rows("18:22").Insert()
$wsBool.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After that insert, the old "trans/VTQaZ/VTQaZ.csv" row (which was row 21)
# is now row 26 -- insert 5 more rows below it for its six entries too.
$wsBool.Rows("27:31").Insert()
$wsBool.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# Carry the formatted-but-empty rows below the data down to rows 33:38.
$wsBool.Range("A32").Copy()
$wsBool.Range("A33:A38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# View state: restore per-sheet selections and the active tab.
# ---------------------------------------------------------------------
$wsInt = $wb.Worksheets.Item("Integer")
$wsInt.Activate()
[void]$wsInt.Range("A13").Select()

$wsBool.Activate()
[void]$wsBool.Range("A32").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
[void]$wsAbout.Range("A1").Select()
